$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Игрушки")

# Fill column F (rows 45-132) with the "Authentication" marker that
# was dropped when the previous merge wasn't finished / file wasn't
# closed before pushing.
$ws.Range("F45:F132").Value = "Authentication"

# Restore the selection/viewport that was active when the fix was made.
$ws.Range("F45:F132").Select()

Write-Host "Done"
